# Daily attendance processing - 2025-12-25 07:57:08
# Rotate the "Recorded By" (column G) list for each data row: move the
# first comma-separated name to the end of the list, leaving single-name
# entries untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($null -eq $value) { continue }

    $text = [string]$value
    if ($text -notmatch ",") { continue }

    $parts = $text -split ",\s*"
    if ($parts.Count -lt 2) { continue }

    $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
    $cell.Value = $rotated
}
